# Roll the 96-quarter "Forecasted Consumption" table forward by one day.
#
# Layout of Sheet1:
#   Row 1            : headers (Timestamp / Forecasted Consumption (MW) / Quarter / Lookup)
#   Rows 2..97        : 96 quarter-hours of the *first* tracked day
#   Rows 98..193      : 96 quarter-hours of the *second* tracked day
#
# The edit drops the oldest day (what used to live in rows 2..97), shifts the
# day that used to live in rows 98..193 up into rows 2..97, and appends a
# brand-new day's worth of forecast numbers into rows 98..193. Column C
# (the 1..96 quarter index) never changes. Column D ("Lookup") is a
# "dd.MM.yyyy" + quarter-number label that rolls forward the same way as the
# timestamp in column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$secondRow = 98
$rowCount = 96
$lastRow = 193

# Brand new forecast (MW) values for the newly appended day (rows 98..193).
$newDayValues = @(
    6270,6240,6210,6180,6140,6120,0,6110,6100,0,
    0,6120,6150,6170,6200,6250,6320,6410,6520,6660,
    6810,6980,7160,7340,7530,7710,7880,8050,8220,8340,
    8430,8520,8590,8600,0,8580,8530,8480,8410,8330,
    8260,8190,8140,8080,8030,8010,8000,0,0,0,
    7990,7960,7930,7910,7900,0,7910,0,7920,7950,
    8000,8060,8120,8180,8270,8340,8410,8500,8580,8620,
    0,0,8600,8580,8560,8520,8480,8410,8330,8230,
    8100,7970,7840,7700,7550,7400,7250,7080,6940,6820,
    6720,6620,6400,6330,6300,6250
)

# --- Step 1: capture the current "second day" (rows 98..193) before it moves ---
$oldSecondA = New-Object 'object[]' $rowCount
$oldSecondB = New-Object 'object[]' $rowCount
$oldSecondD = New-Object 'object[]' $rowCount
for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $secondRow + $i
    # Value2 returns the raw double/string instead of auto-coercing
    # date-formatted numbers into a (lower precision) DateTime.
    $oldSecondA[$i] = $ws.Range("A$r").Value2()
    $oldSecondB[$i] = $ws.Range("B$r").Value2()
    $oldSecondD[$i] = $ws.Range("D$r").Value2()
}

# --- Step 2: move that captured day into the first-day rows (2..97) ---
$firstArrAB = New-Object 'object[,]' $rowCount,2
$firstArrD = New-Object 'object[,]' $rowCount,1
for ($i = 0; $i -lt $rowCount; $i++) {
    $firstArrAB[$i,0] = [double]$oldSecondA[$i]
    $firstArrAB[$i,1] = [double]$oldSecondB[$i]
    $firstArrD[$i,0] = [string]$oldSecondD[$i]
}
$firstLastRow = $firstRow + $rowCount - 1
$ws.Range("A" + $firstRow + ":B" + $firstLastRow).Value = $firstArrAB
$ws.Range("D" + $firstRow + ":D" + $firstLastRow).Value = $firstArrD

# --- Step 3: build the new second day (rows 98..193): date +1, brand-new
#             consumption figures, and the Lookup label rolled to the new date ---
$secondArrAB = New-Object 'object[,]' $rowCount,2
$secondArrD = New-Object 'object[,]' $rowCount,1
for ($i = 0; $i -lt $rowCount; $i++) {
    $secondArrAB[$i,0] = [double]$oldSecondA[$i] + 1
    $secondArrAB[$i,1] = [double]$newDayValues[$i]
    $oldLabel = [string]$oldSecondD[$i]
    $secondArrD[$i,0] = $oldLabel.Replace("17.02.2026", "18.02.2026")
}
$secondLastRow = $secondRow + $rowCount - 1
$ws.Range("A" + $secondRow + ":B" + $secondLastRow).Value = $secondArrAB
$ws.Range("D" + $secondRow + ":D" + $secondLastRow).Value = $secondArrD
